$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182674288749695
$ws.Range("B1").Value = 2.407591342926025
$ws.Range("C1").Value = 3.83289098739624
$ws.Range("D1").Value = 2.139377593994141
$ws.Range("E1").Value = 1.200313806533813
